$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 3: fill in the rest of the "nostril" test entry ---
$ws.Range("B3").Value = "Error"
$ws.Range("C3").Value = "Traceback (most recent call last):" + [char]10 + `
  'File "C:\Users\Shannon\Documents\GitHub\DVS-Python\eyeDetection.py", line 99, in <module>' + [char]10 + `
  '    image = DetectRedEyes(img, faceCascade, eyeCascade)' + [char]10 + `
  'File "C:\Users\Shannon\Documents\GitHub\DVS-Python\eyeDetection.py", line 71, in DetectRedEyes' + [char]10 + `
  '    cv.SetImageROI(image, (pt1[0],' + [char]10 + `
  "UnboundLocalError: local variable 'pt1' referenced before assignment"
$ws.Range("E3").Value = "Shannon Harris"

# D3 needs the same date-formatted style already used by D2 (numFmt 14);
# copying the cell first picks up that style, then we overwrite the value.
$ws.Range("D2").Copy($ws.Range("D3"))
$ws.Range("D3").Value = 41389

# This row now holds the full (wrapped) traceback, so it needs to grow.
$ws.Rows.Item(3).RowHeight = 210

# --- Row 20: two more photos were tested; record a successful run ---
$ws.Range("B20").Value = "Success"
$ws.Range("C20").Value = "None"
$ws.Range("D20").Value = 41389
$ws.Range("E20").Value = "Shannon Harris"

$ws.Range("D2").Copy($ws.Range("D20"))
$ws.Range("D20").Value = 41389

# --- Scroll the view down a bit and land the selection on D4 ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 3
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D4").Select()
